$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The original sheet had 3 columns: A=fecha, B=tarea, C=responsables.
# The new layout splits the date column into "fecha incio" / "fecha fin",
# pushing the old "tarea" / "responsables" columns one slot to the right.
# Insert a fresh column at B so old-B -> C and old-C -> D.
$ws.Columns.Item(2).Insert()

# --- Header row (row 4) ---
$ws.Range("A4").Value = "fecha incio"
$ws.Range("B4").Value = "fecha fin"
# C4 ("tarea") and D4 ("responsables") already shifted into place.

# --- Row 5: Analisis de requerimiento ---
$ws.Range("B5").Value = 45518
$ws.Range("B5").NumberFormat = "d-mmm"
$ws.Range("C5").Value = "Analisis de requerimiento -FODA"
$ws.Range("D5").Value = "todos"

# --- Row 6: Conseguir un cliente ---
$ws.Range("A6").Value = 45518
$ws.Range("A6").NumberFormat = "d-mmm"
$ws.Range("B6").Value = 45596
$ws.Range("B6").NumberFormat = "d-mmm"
$ws.Range("C6").Value = "Conseguir un cliente"
$ws.Range("D6").Value = "Engers"

# --- Row 7: Buscar productos relacionados directamente ---
$ws.Range("B7").Value = 45596
$ws.Range("B7").NumberFormat = "d-mmm"
$ws.Range("C7").Value = "Buscar productos relacionados directamente"
$ws.Range("D7").Value = "Lima, Gomez"

# --- Row 8: Buscar productos relacionados indirectamente ---
$ws.Range("B8").Value = 45596
$ws.Range("B8").NumberFormat = "d-mmm"
$ws.Range("C8").Value = "Buscar productos relacionados indirectamente"
$ws.Range("D8").Value = "Lima, Gomez"

# --- Row 9: Establecer Vision ---
$ws.Range("B9").Value = 45561
$ws.Range("B9").NumberFormat = "d-mmm"
$ws.Range("C9").Value = "Establecer Vision"
$ws.Range("D9").Value = "Lima"

# --- Column widths ---
# Columns A & B share the original "fecha" column width (~19.71 chars).
$ws.Columns.Item(1).ColumnWidth = 18.8
$ws.Columns.Item(2).ColumnWidth = 18.8
# Column C (tarea) grew to fit the longest task description (~42.71 chars).
$ws.Columns.Item(3).ColumnWidth = 41.8
# Column D keeps the original "responsables" width untouched (inherited from old column C).

# --- Selection matches the author's last-edited cell ---
$ws.Range("B9").Select() | Out-Null
